# Adapt column header formatting to respective input file names (#7)
#   * rename "<name>_old" headers to "<name>_FV2310"
#   * rename "<name>_new" headers to "<name>_FV2404"
#   * wrap the data range in an Excel Table ("Table1")
#   * freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Turn the used range into a proper Excel Table (ListObject) with the
# renamed headers, matching the ref "A1:U59" used by the sheet's data.
$dataRange = $ws.Range("A1:U59")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Freeze the header row (pane split after row 1).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
